$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct a couple of attraction name typos/formatting in column A
$ws.Range("A11").Value = "ET Adventure"
$ws.Range("A10").Value = "TRANSFORMERS - The Ride-3D"

# Update the active selection to A10
$ws.Range("A10").Select()
